$d = $word.ActiveDocument
$d.Content.Find.Execute("or any other metal would attack the mold", $true, $false, $false, $false, $false,
                         $true, 1, $false, "or any other metal would attach to the mold", 2)
